# Insert a new data row before the current row 742 ("Primera" quality,
# Region del Maule, date 2022-03-02) which pushes all rows 742:839 down to
# 743:840. Then fill in the newly inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("742:742").Insert()

$ws.Cells.Item(742, 1).Value = 10
$ws.Cells.Item(742, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(742, 3).Value = "La Araucanía"
$ws.Cells.Item(742, 4).Value = 45131
$ws.Cells.Item(742, 5).Value = 9
$ws.Cells.Item(742, 6).Value = 100112043
$ws.Cells.Item(742, 7).Value = "Pepino ensalada"
$ws.Cells.Item(742, 8).Value = "Sin especificar"
$ws.Cells.Item(742, 9).Value = "Primera"
$ws.Cells.Item(742, 10).Value = 235
$ws.Cells.Item(742, 11).Value = 13000
$ws.Cells.Item(742, 12).Value = 14000
$ws.Cells.Item(742, 13).Value = 13468
$ws.Cells.Item(742, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(742, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(742, 16).Value = 269
$ws.Cells.Item(742, 17).Value = 50
$ws.Cells.Item(742, 18).Value = "Hortaliza"
